$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (ownTeam, oppTeam) before the old column D (batsman),
# shifting batsman..sr from D:I to F:K
$ws.Range("D1:E1").EntireColumn.Insert()

# Force the whole target range to Text format so numeric-looking values
# (runs, balls, 4s, 6s, sr) are stored as text, matching the source data
$ws.Range("A1:K12").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Data rows
# Row 2
$ws.Range("A2").Value = " Dubai (DSC)"
$ws.Range("B2").Value = " October 27 2020"
$ws.Range("C2").Value = "Sunrisers won by 88 runs"
$ws.Range("D2").Value = "Sunrisers Hyderabad"
$ws.Range("E2").Value = "Delhi Capitals"
$ws.Range("F2").Value = "Kane Williamson "
$ws.Range("G2").Value = "11"
$ws.Range("H2").Value = "10"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "110.00"

# Row 3
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 13 2020"
$ws.Range("C3").Value = "Super Kings won by 20 runs"
$ws.Range("D3").Value = "Sunrisers Hyderabad"
$ws.Range("E3").Value = "Chennai Super Kings"
$ws.Range("F3").Value = "Kane Williamson "
$ws.Range("G3").Value = "57"
$ws.Range("H3").Value = "39"
$ws.Range("I3").Value = "7"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "146.15"

# Row 4
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 18 2020"
$ws.Range("C4").Value = "Match tied (KKR won the one-over eliminator)"
$ws.Range("D4").Value = "Sunrisers Hyderabad"
$ws.Range("E4").Value = "Kolkata Knight Riders"
$ws.Range("F4").Value = "Kane Williamson "
$ws.Range("G4").Value = "29"
$ws.Range("H4").Value = "19"
$ws.Range("I4").Value = "4"
$ws.Range("J4").Value = "1"
$ws.Range("K4").Value = "152.63"

# Row 5
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " November 06 2020"
$ws.Range("C5").Value = "Sunrisers won by 6 wickets (with 2 balls remaining)"
$ws.Range("D5").Value = "Sunrisers Hyderabad"
$ws.Range("E5").Value = "Royal Challengers Bangalore"
$ws.Range("F5").Value = "Kane Williamson "
$ws.Range("G5").Value = "50"
$ws.Range("H5").Value = "44"
$ws.Range("I5").Value = "2"
$ws.Range("J5").Value = "2"
$ws.Range("K5").Value = "113.63"

# Row 6
$ws.Range("A6").Value = " Sharjah"
$ws.Range("B6").Value = " October 31 2020"
$ws.Range("C6").Value = "Sunrisers won by 5 wickets (with 35 balls remaining)"
$ws.Range("D6").Value = "Sunrisers Hyderabad"
$ws.Range("E6").Value = "Royal Challengers Bangalore"
$ws.Range("F6").Value = "Kane Williamson "
$ws.Range("G6").Value = "8"
$ws.Range("H6").Value = "14"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "57.14"

# Row 7
$ws.Range("A7").Value = " Abu Dhabi"
$ws.Range("B7").Value = " November 08 2020"
$ws.Range("C7").Value = "Capitals won by 17 runs"
$ws.Range("D7").Value = "Sunrisers Hyderabad"
$ws.Range("E7").Value = "Delhi Capitals"
$ws.Range("F7").Value = "Kane Williamson "
$ws.Range("G7").Value = "67"
$ws.Range("H7").Value = "45"
$ws.Range("I7").Value = "5"
$ws.Range("J7").Value = "4"
$ws.Range("K7").Value = "148.88"

# Row 8
$ws.Range("A8").Value = " Dubai (DSC)"
$ws.Range("B8").Value = " October 11 2020"
$ws.Range("C8").Value = "Royals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D8").Value = "Sunrisers Hyderabad"
$ws.Range("E8").Value = "Rajasthan Royals"
$ws.Range("F8").Value = "Kane Williamson "
$ws.Range("G8").Value = "22"
$ws.Range("H8").Value = "12"
$ws.Range("I8").Value = "0"
$ws.Range("J8").Value = "2"
$ws.Range("K8").Value = "183.33"

# Row 9
$ws.Range("A9").Value = " Dubai (DSC)"
$ws.Range("B9").Value = " October 08 2020"
$ws.Range("C9").Value = "Sunrisers won by 69 runs"
$ws.Range("D9").Value = "Sunrisers Hyderabad"
$ws.Range("E9").Value = "Kings XI Punjab"
$ws.Range("F9").Value = "Kane Williamson "
$ws.Range("G9").Value = "20"
$ws.Range("H9").Value = "10"
$ws.Range("I9").Value = "1"
$ws.Range("J9").Value = "1"
$ws.Range("K9").Value = "200.00"

# Row 10
$ws.Range("A10").Value = " Dubai (DSC)"
$ws.Range("B10").Value = " October 02 2020"
$ws.Range("C10").Value = "Sunrisers won by 7 runs"
$ws.Range("D10").Value = "Sunrisers Hyderabad"
$ws.Range("E10").Value = "Chennai Super Kings"
$ws.Range("F10").Value = "Kane Williamson "
$ws.Range("G10").Value = "9"
$ws.Range("H10").Value = "13"
$ws.Range("I10").Value = "1"
$ws.Range("J10").Value = "0"
$ws.Range("K10").Value = "69.23"

# Row 11
$ws.Range("A11").Value = " Sharjah"
$ws.Range("B11").Value = " October 04 2020"
$ws.Range("C11").Value = "Mumbai won by 34 runs"
$ws.Range("D11").Value = "Sunrisers Hyderabad"
$ws.Range("E11").Value = "Mumbai Indians"
$ws.Range("F11").Value = "Kane Williamson "
$ws.Range("G11").Value = "3"
$ws.Range("H11").Value = "5"
$ws.Range("I11").Value = "0"
$ws.Range("J11").Value = "0"
$ws.Range("K11").Value = "60.00"

# Row 12
$ws.Range("A12").Value = " Abu Dhabi"
$ws.Range("B12").Value = " September 29 2020"
$ws.Range("C12").Value = "Sunrisers won by 15 runs"
$ws.Range("D12").Value = "Sunrisers Hyderabad"
$ws.Range("E12").Value = "Delhi Capitals"
$ws.Range("F12").Value = "Kane Williamson "
$ws.Range("G12").Value = "41"
$ws.Range("H12").Value = "26"
$ws.Range("I12").Value = "5"
$ws.Range("J12").Value = "0"
$ws.Range("K12").Value = "157.69"

